{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block (and\n// the blank paragraph right before it) that used to follow the last\n// \"Requisitos\" entry (\"LOB1215: Recursos Energ\u00e9ticos (Requisito fraco)\").\n// The blank paragraph that originally sat *after* the footer (right before\n// the page-break paragraph) is left untouched.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph the footer block follows.\nconst anchorText = \"LOB1215: Recursos Energ\u00e9ticos (Requisito fraco)\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// The three paragraphs immediately following the anchor are the ones being\n// removed: an empty paragraph, the \"Ver no Jupiter ...\" line, and the\n// \"\u00a9 2020 ...\" copyright line.\nconst expectedTexts = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst toDelete = [];\nfor (let k = 0; k < expectedTexts.length; k++) {\n  const para = items[anchorIndex + 1 + k];\n  if (!para || para.text.trim() !== expectedTexts[k]) {\n    throw new Error(\n      \"Unexpected document structure near footer block; expected '\" +\n        expectedTexts[k] +\n        \"' but found '\" +\n        (para ? para.text : \"<none>\") +\n        \"'\"\n    );\n  }\n  toDelete.push(para);\n}\n\n// Delete from the last one back to the first so earlier references stay valid.\nfor (let k = toDelete.length - 1; k >= 0; k--) {\n  toDelete[k].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" footer block (and\n# the blank paragraph right before it) that used to follow the last\n# \"Requisitos\" entry (\"LOB1215: Recursos Energeticos (Requisito fraco)\").\n# The blank paragraph that originally sat *after* the footer (right before\n# the page-break paragraph) is left untouched.\n\nfunction Normalize-ParaText($t) {\n    if ($null -eq $t) { return \"\" }\n    return $t.Trim([char]13, [char]7, [char]10, ' ')\n}\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph. Match on the ASCII-safe prefix only, since\n# accented characters can round-trip lossily through this text channel.\n$count = $d.Paragraphs.Count\n$anchor = $null\nfor ($i = 1; $i -le $count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ((Normalize-ParaText $candidate.Range.Text) -like \"LOB1215:*\") {\n        $anchor = $candidate\n        break\n    }\n}\nif ($null -eq $anchor) {\n    throw \"Could not find anchor paragraph starting with 'LOB1215:'\"\n}\n\n# The three paragraphs directly following the anchor make up the footer\n# block being removed: an empty paragraph, the \"Ver no Jupiter ...\" line,\n# and the \"(c) 2020 ...\" copyright line.\n$p1 = $anchor.Next()\n$p2 = $p1.Next()\n$p3 = $p2.Next()\n\nif ((Normalize-ParaText $p1.Range.Text) -ne \"\") {\n    throw \"Unexpected content where the blank paragraph was expected: '$($p1.Range.Text)'\"\n}\nif ((Normalize-ParaText $p2.Range.Text) -notlike \"Ver no Jupiter*\") {\n    throw \"Unexpected content where the 'Ver no Jupiter ...' paragraph was expected: '$($p2.Range.Text)'\"\n}\nif ((Normalize-ParaText $p3.Range.Text) -notlike \"*Contact: luizeleno@usp.br*\") {\n    throw \"Unexpected content where the copyright paragraph was expected: '$($p3.Range.Text)'\"\n}\n\n# Delete starting from the last paragraph so the earlier references stay valid.\n$p3.Range.Delete()\n$p2.Range.Delete()\n$p1.Range.Delete()\n"}
